# Updates the cryptocurrency price/volume snapshot pulled in by the
# "Updated cryptos list" GitHub Actions workflow.
#
# Most rows only get refreshed Price (column D) / Volume(1h) (column E)
# figures. Rows 47-48 additionally swap rank between "Stellar" and
# "FirstDigitalUSD" (Name/Link/Price all change).
#
# Price values such as "1.00" / "0.997" / "0.0404" are textual (the sheet
# keeps them as strings, some even use "." as a thousands separator, e.g.
# "66.648.85"), so for any replacement that Excel would otherwise parse as
# a genuine number we force the cell to Text format first to keep the
# exact digits/trailing zeros instead of Excel silently re-typing it as a
# number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '66.648.85'
$ws.Range("E2").Value = '  +1.88%  '
# Row 3
$ws.Range("D3").Value = '3.264.50'
$ws.Range("E3").Value = '  -0.68%  '
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.997'
$ws.Range("E4").Value = '  -0.26%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '568.47'
$ws.Range("E5").Value = '  -1.44%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.82'
$ws.Range("E6").Value = '  -4.16%  '
# Row 7
$ws.Range("E7").Value = '  -0.06%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.581'
$ws.Range("E8").Value = '  +2.01%  '
# Row 9
$ws.Range("D9").Value = '3.258.56'
$ws.Range("E9").Value = '  -0.74%  '
# Row 10
$ws.Range("E10").Value = '  -1.33%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.568'
$ws.Range("E11").Value = '  -0.50%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.13'
$ws.Range("E12").Value = '  -2.79%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000267'
$ws.Range("E13").Value = '  +1.79%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '691.78'
$ws.Range("E14").Value = '  +10.29%  '
# Row 15
$ws.Range("D15").Value = '3.780.55'
$ws.Range("E15").Value = '  -0.81%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.27'
$ws.Range("E16").Value = '  -1.30%  '
# Row 17
$ws.Range("D17").Value = '66.670.88'
$ws.Range("E17").Value = '  +1.72%  '
# Row 18
$ws.Range("E18").Value = '  +0.95%  '
# Row 19
$ws.Range("D19").Value = '3.265.46'
$ws.Range("E19").Value = '  -0.69%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.23'
$ws.Range("E20").Value = '  -2.07%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.68'
$ws.Range("E21").Value = '  -2.01%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.883'
$ws.Range("E22").Value = '  -0.19%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '16.83'
$ws.Range("E23").Value = '  -6.19%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.11'
$ws.Range("E24").Value = '  +3.18%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '97.67'
$ws.Range("E25").Value = '  -2.05%  '
# Row 26
$ws.Range("E26").Value = '  -2.43%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.69'
$ws.Range("E27").Value = '  -1.21%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.24'
$ws.Range("E28").Value = '  -0.91%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '32.65'
$ws.Range("E29").Value = '  +6.84%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.35'
$ws.Range("E30").Value = '  +0.39%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.67'
$ws.Range("E31").Value = '  +3.45%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '574.89'
$ws.Range("E32").Value = '  +1.10%  '
# Row 33
$ws.Range("D33").Value = '3.857.82'
$ws.Range("E33").Value = '  +0.57%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '10.73'
$ws.Range("E34").Value = '  -0.93%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.103'
$ws.Range("E35").Value = '  -0.60%  '
# Row 36
$ws.Range("E36").Value = '  -0.01%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '55.17'
$ws.Range("E37").Value = '  -0.51%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.26'
$ws.Range("E38").Value = '  -11.18%  '
# Row 39
$ws.Range("E39").Value = '  +1.98%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.59'
$ws.Range("E40").Value = '  +0.41%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.34'
$ws.Range("E41").Value = '  -1.63%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '31.54'
$ws.Range("E42").Value = '  -2.48%  '
# Row 43
$ws.Range("D43").Value = '0.0₃0662'
$ws.Range("E43").Value = '  -2.26%  '
# Row 44
$ws.Range("E44").Value = '  -2.77%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.325'
$ws.Range("E45").Value = '  -1.42%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0404'
$ws.Range("E46").Value = '  +0.13%  '
# Row 47
$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.127'
$ws.Range("E47").Value = '  +0.48%  '
# Row 48
$ws.Range("B48").Value = 'FirstDigitalUSD'
$ws.Range("C48").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  +0.10%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.36'
$ws.Range("E49").Value = '  +7.74%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.50'
$ws.Range("E50").Value = '  -0.24%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '129.11'
$ws.Range("E51").Value = '  +0.02%  '
